$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (F3:G3) - lowercase text + yellow fill
$ws.Range("F3").Value = "romaji"
$ws.Range("G3").Value = "english"
$ws.Range("F3:G3").Interior.Color = 65535

# Data rows F4:G48 - re-sorted + newly added vocabulary
$data = New-Object 'object[,]' 45,2
$data[0,0] = "ani mo urusai desu"
$data[0,1] = "my older brother is also noisy"
$data[1,0] = "basuketto booru o shimasu ka"
$data[1,1] = "do you play basketball"
$data[2,0] = "basukettoboru o shimasuka"
$data[2,1] = "do you play basketball"
$data[3,0] = "chichi mo o share desu"
$data[3,1] = "my dad is akso fashionable"
$data[4,0] = "chichi wa toronto shusshin desu"
$data[4,1] = "my dad is from toronto"
$data[5,0] = "hahamoisogashiidesu"
$data[5,1] = "my mom is also busy"
$data[6,0] = "hai terebi o mimasu"
$data[6,1] = "yes i watch tv"
$data[7,0] = "hajimemashite"
$data[7,1] = "nice to meet you "
$data[8,0] = "hana san wa totemo urusai desu"
$data[8,1] = "hana is very noisy"
$data[9,0] = "hon o yo mimasuka"
$data[9,1] = "do you read books"
$data[10,0] = "iie jazu o kikimasu"
$data[10,1] = "no i listen to jazz"
$data[11,0] = "imouto wa akarui hito desu"
$data[11,1] = "my younger sister is a cheerful person"
$data[12,0] = "imouto wa daigakusei desu"
$data[12,1] = "my younger sister is a college student"
$data[13,0] = "jpoppu okikimasuka"
$data[13,1] = "do you listen to jpop"
$data[14,0] = "ken san wa nyuu you ku sshu shin desu"
$data[14,1] = "ken is from new york"
$data[15,0] = "korewa anime janais desu"
$data[15,1] = "this is not anime"
$data[16,0] = "manga o yomimaska"
$data[16,1] = "do you read mangas"
$data[17,0] = "musuko wa yumei desu"
$data[17,1] = "my son is famous"
$data[18,0] = "musume wa bengoshi desu"
$data[18,1] = "my daughter is a lawyer"
$data[19,0] = "musume wa gakusei desu"
$data[19,1] = "my daughter is a student"
$data[20,0] = "nakayama san wa yoku yaku o shimasu"
$data[20,1] = "Ms nakayama often plays baseball"
$data[21,0] = "naomisannokazokumochisai desu"
$data[21,1] = "naomis family is also small"
$data[22,0] = "nihon ni sundeimasu"
$data[22,1] = "i live in japan"
$data[23,0] = "nyuu yoo ku wa chisai toshi janaii desu"
$data[23,1] = "new yoors is not a small city"
$data[24,0] = "nyuu yoo ku wa totemo ooki desu"
$data[24,1] = "new york  is very big"
$data[25,0] = "ohayou gozaimasu"
$data[25,1] = "good morning"
$data[26,0] = "otoutowa atamagaiidesu"
$data[26,1] = "my youger brother is smart"
$data[27,0] = "otoutowa nyuuyooku shusshin desu"
$data[27,1] = "my husband is from new york"
$data[28,0] = "rokku o kikimasu"
$data[28,1] = "i listen to rock"
$data[29,0] = "sonohito wa kaishain desu ka"
$data[29,1] = "is that person an office worker"
$data[30,0] = "terebitoeiga"
$data[30,1] = "tv and movies"
$data[31,0] = "tokidoki hon o yomimasuka"
$data[31,1] = "do you sometimes read book"
$data[32,0] = "tokidoki nyusu o yomimasu"
$data[32,1] = "i sometimes read news"
$data[33,0] = "tokidoki sakaa o shimasu"
$data[33,1] = "i sometimes play soccer"
$data[34,0] = "toronto ni sundeimasuka"
$data[34,1] = "do you live in toronto"
$data[35,0] = "toukyouni sundeimasuka"
$data[35,1] = "do you live in tokyo"
$data[36,0] = "tsuma wa totemo oshara desu"
$data[36,1] = "my wife is very fashionable"
$data[37,0] = "watashimo kaishain desu"
$data[37,1] = "i am also an office worker"
$data[38,0] = "yakyu u o shimasu ka"
$data[38,1] = "do you play baseball"
$data[39,0] = "yoku hon o yomimasu"
$data[39,1] = "i often read books"
$data[40,0] = "yoku ongaku o kikimasu ka"
$data[40,1] = "do you often listen to music"
$data[41,0] = "ongaku o kikimasuka"
$data[41,1] = "do you listen to music"
$data[42,0] = "yakyu u o shimasu"
$data[42,1] = "i play baseball"
$data[43,0] = "yoku eiga o mimasuka"
$data[43,1] = "do you often watch movies"
$data[44,0] = "naomisanwatokidokianimeomimasu"
$data[44,1] = ""
$ws.Range("F4:G48").Value = $data

# Restore selection to F24 (per saved workbook view state)
$ws.Range("F24").Select()